{"js": "// The document has several \"<id>...</id>\" markers that were typed as three\n// separate runs: \"<id>\", \"<the-id-value>\", \"</id>\". This collapses each of\n// those triples into a single run (keeping the \"<id>\"/\"</id>\" run's\n// formatting) whose text is the full \"<id>value</id>\" string, matching how\n// the surrounding markers (e.g. \"<page>...</page>\") are already encoded.\nconst idValues = [\"p052r_1\", \"p052v_1\", \"p052v_2\"];\n\nconst body = context.document.body;\n\nfor (const value of idValues) {\n  const marker = `<id>${value}</id>`;\n  const results = body.search(marker, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    // Replacing the whole matched range with the identical text merges the\n    // three runs into a single run, inheriting the formatting of the first\n    // (the \"<id>\"/\"</id>\" Courier-New run), exactly like the native edit.\n    results.items[i].insertText(marker, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document encodes several \"<id>...</id>\" markers as three separate\n# runs: \"<id>\", \"<the-id-value>\", \"</id>\". Collapse each triple into a\n# single run (keeping the \"<id>\"/\"</id>\" run's Courier-New formatting)\n# whose text is the full \"<id>value</id>\" string - matching how sibling\n# markers such as \"<page>...</page>\" are already encoded as one run.\n\n$d = $word.ActiveDocument\n$idValues = @(\"p052r_1\", \"p052v_1\", \"p052v_2\")\n\nforeach ($val in $idValues) {\n    $marker = \"<id>\" + $val + \"</id>\"\n\n    # $rng gets mutated in place to the matched range by Find.Execute().\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $marker\n    $find.MatchCase = $true\n    $found = $find.Execute()\n\n    if ($found) {\n        $markerStart = $rng.Start\n        $openLen = 4                 # length of \"<id>\"\n        $restStart = $markerStart + $openLen\n\n        # Everything after \"<id>\" (the id value + \"</id>\") currently lives\n        # in two more runs; delete that text and retype it onto the end of\n        # the \"<id>\" run so the three runs become one.\n        $restLen = $val.Length + 5   # id value + length of \"</id>\"\n        $rRest = $d.Range($restStart, $restStart + $restLen)\n        $rRest.Delete()\n\n        $rFirst = $d.Range($markerStart, $restStart)\n        $rFirst.InsertAfter($val + \"</id>\")\n    }\n}\n"}
